$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column in H, mirroring the formatting used by the
# existing header cells (copy format from G1 so the header shares the
# same style as the rest of row 1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"

# Fill in the "Save" values for each data row (2-12)
$saveValues = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
